$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Swap the two adjacent empty bookmarks so that the one named
#    "_Hlk500579323" ends up first (id 1) and "_Hlk500579466" second
#    (id 2) -- matching the target XML ordering.
# ------------------------------------------------------------------
$bmA = $d.Bookmarks.Item("_Hlk500579466")
$rangeA = $bmA.Range
$bmB = $d.Bookmarks.Item("_Hlk500579323")
$rangeB = $bmB.Range
$bmA.Delete()
$bmB.Delete()
$d.Bookmarks.Add("_Hlk500579466", $rangeA)
$d.Bookmarks.Add("_Hlk500579323", $rangeB)

# ------------------------------------------------------------------
# 2) Fix the SGF2022 year typo: "Apr 10 - Apr 13, 2021?" should read
#    "Apr 10 - Apr 13, 2022?"
# ------------------------------------------------------------------
$c = $d.Content
$found = $c.Find.Execute("Apr 10 – Apr 13, 2021?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $c.End
$oneChar = $d.Range($endPos - 2, $endPos - 1)
$oneChar.Delete()
$insertPos = $d.Range($endPos - 2, $endPos - 2)
$insertPos.InsertAfter("2")
$touch = $d.Range($endPos - 2, $endPos - 1)
$touch.Bold = 1
$touch.Bold = 0

# ------------------------------------------------------------------
# 3) Adjust the "Table 2: Capabilities by Conference Content topic"
#    column widths by 1 dxa: col1 1543->1542, col2 1872->1873.
# ------------------------------------------------------------------
$t = $d.Tables.Item(3)
$col1 = $t.Columns.Item(1)
$col1.Width = [double](1542.0 / 20.0)
$col2 = $t.Columns.Item(2)
$col2.Width = [double](1873.0 / 20.0)
